$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.554.52'
$ws.Range("E2").Value = '  +1.57%  '
$ws.Range("D3").Value = '3.166.21'
$ws.Range("E3").Value = '  +1.32%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.158.52'
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("E9").Value = '  -0.21%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.162'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("E11").Value = '  +3.99%  '
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.99%  '
$ws.Range("D15").Value = '3.690.39'
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("E16").Value = '  -1.04%  '
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '64.327.94'
$ws.Range("E18").Value = '  +1.33%  '
$ws.Range("D19").Value = '3.161.69'
$ws.Range("E19").Value = '  +1.16%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.93'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.43%  '
$ws.Range("E22").Value = '  +0.53%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.17'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("B25").Value = 'Fetch.AI'
$ws.Range("C25").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.36'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.65'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.36%  '
$ws.Range("E27").Value = '  +0.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.79'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +10.03%  '
$ws.Range("E29").Value = '  +2.23%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.74'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.63%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.70'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("E34").Value = '  +1.37%  '
$ws.Range("D35").Value = '0.0₃0847'
$ws.Range("E35").Value = '  -2.50%  '
$ws.Range("E36").Value = '  +2.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.29'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.33'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '474.21'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.56%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '9.29'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.297'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.89%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0376'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").Value = '2.933.83'
$ws.Range("E45").Value = '  +1.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +13.44%  '
$ws.Range("E47").Value = '  -1.69%  '
$ws.Range("E48").Value = '  +3.78%  '
$ws.Range("E50").Value = '  +3.67%  '
$ws.Range("E51").Value = '  -0.07%  '
